# Updated lai, cover and node number data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CottonObserved")

# Fill in newly observed data for row 18 (Cotton.Leaf.NodeNumber, TotalNumber, LAI,
# LAIError, CoverGreen, CoverGreenError)
$ws.Range("D18").Value = 19.649999999999999
$ws.Range("E18").Value = 66.849999999999994
$ws.Range("F18").Value = 5.312378794687266
$ws.Range("G18").Value = 0.32232021196016941
$ws.Range("H18").Value = 0.98990384615384619
$ws.Range("H18").NumberFormat = "0.00"
$ws.Range("I18").Value = 0.063817475023250456

# Update view: move the active selection to reflect the newly entered data
$ws.Range("L18").Select()
